$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N. This pushes the existing
# "Late" / "heading" / "Outstanding" columns one slot to the right
# (N->O, O->P, P->Q) and leaves the new N column blank.
$null = $ws.Columns("N").Insert()

# The newly inserted column inherits the width of its left neighbour (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth()

# "Repayment schedule" becomes the active sheet/tab, with K13 selected
# (previously "Edit Repayment Schedule" was the active tab).
$null = $ws.Activate()
$null = $ws.Range("K13").Select()
